# Auto-generated edit script applying scheduled market-data refresh values
# to the Leve profit-calculation columns (H:N) across the class worksheets.
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 181.38889
$ws.Range("I42").Value = 63.916668
$ws.Range("J42").Value = 416.33334
$ws.Range("K42").Value = 191.750004
$ws.Range("L42").Value = 1249.00002
$ws.Range("M42").Value = 38.24999600000001
$ws.Range("N42").Value = -1709.00002
$ws.Range("H129").Value = 1137.7925
$ws.Range("J129").Value = 1197.0613
$ws.Range("L129").Value = 3591.1839
$ws.Range("N129").Value = -13591.1839

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1933.3334
$ws.Range("I5").Value = 1933.3334
$ws.Range("K5").Value = 1933.3334
$ws.Range("M5").Value = -1821.3334
$ws.Range("H88").Value = 2797.1428
$ws.Range("I88").Value = 2648
$ws.Range("J88").Value = 2856.8
$ws.Range("K88").Value = 2648
$ws.Range("L88").Value = 2856.8
$ws.Range("M88").Value = -2242
$ws.Range("N88").Value = -3668.8
$ws.Range("H91").Value = 2797.1428
$ws.Range("I91").Value = 2648
$ws.Range("J91").Value = 2856.8
$ws.Range("K91").Value = 2648
$ws.Range("L91").Value = 2856.8
$ws.Range("M91").Value = -1244
$ws.Range("N91").Value = -5664.8

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1933.3334
$ws.Range("I4").Value = 1933.3334
$ws.Range("K4").Value = 1933.3334
$ws.Range("M4").Value = -1818.3334

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 40000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 40000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 40000
$ws.Range("N32").Value = -40632
$ws.Range("M32").ClearContents()
$ws.Range("H134").Value = 2062.1875
$ws.Range("I134").Value = 1750.2609
$ws.Range("J134").Value = 2859.3333
$ws.Range("K134").Value = 5250.7827
$ws.Range("L134").Value = 8577.999899999999
$ws.Range("M134").Value = -2715.7827
$ws.Range("N134").Value = -13647.9999
$ws.Range("H137").Value = 34189.832
$ws.Range("J137").Value = 34189.832
$ws.Range("L137").Value = 34189.832
$ws.Range("N137").Value = -44389.832

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9675.675999999999
$ws.Range("J87").Value = 9942.857
$ws.Range("L87").Value = 29828.571
$ws.Range("N87").Value = -32324.571
$ws.Range("H90").Value = 9675.675999999999
$ws.Range("J90").Value = 9942.857
$ws.Range("L90").Value = 89485.713
$ws.Range("N90").Value = -101965.713
$ws.Range("H103").Value = 3122.2856
$ws.Range("J103").Value = 3311.2
$ws.Range("L103").Value = 9933.599999999999
$ws.Range("N103").Value = -11691.6
$ws.Range("H117").Value = 83938.5
$ws.Range("I117").Value = 533
$ws.Range("J117").Value = 167344
$ws.Range("K117").Value = 1599
$ws.Range("L117").Value = 502032
$ws.Range("M117").Value = 1843
$ws.Range("N117").Value = -508916
$ws.Range("H120").Value = 14375
$ws.Range("I120").Value = 8750
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 26250
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -21412
$ws.Range("N120").Value = -69676
$ws.Range("H121").Value = 33259.59
$ws.Range("I121").Value = 947
$ws.Range("J121").Value = 79420.42999999999
$ws.Range("K121").Value = 2841
$ws.Range("L121").Value = 238261.29
$ws.Range("M121").Value = -1531
$ws.Range("N121").Value = -240881.29
$ws.Range("H122").Value = 1415.3846
$ws.Range("I122").Value = 1433.3334
$ws.Range("J122").Value = 1413.0435
$ws.Range("K122").Value = 12900.0006
$ws.Range("L122").Value = 12717.3915
$ws.Range("M122").Value = -10450.0006
$ws.Range("N122").Value = -17617.3915
$ws.Range("H123").Value = 7000
$ws.Range("I123").Value = 4000
$ws.Range("J123").Value = 10000
$ws.Range("K123").Value = 12000
$ws.Range("L123").Value = 30000
$ws.Range("M123").Value = -9550
$ws.Range("N123").Value = -34900
$ws.Range("H124").Value = 12440
$ws.Range("I124").Value = 1030
$ws.Range("J124").Value = 14070
$ws.Range("K124").Value = 3090
$ws.Range("L124").Value = 42210
$ws.Range("M124").Value = 1820
$ws.Range("N124").Value = -52030
$ws.Range("H125").Value = 3523.75
$ws.Range("I125").Value = 300
$ws.Range("J125").Value = 3984.2856
$ws.Range("K125").Value = 900
$ws.Range("L125").Value = 11952.8568
$ws.Range("M125").Value = 4020
$ws.Range("N125").Value = -21792.8568
$ws.Range("H126").Value = 2343.3333
$ws.Range("I126").Value = 2343.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7029.999899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2089.999899999999
$ws.Range("H127").Value = 818.8333
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 818.8333
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 2456.4999
$ws.Range("N127").Value = -12376.4999
$ws.Range("H128").Value = 170000
$ws.Range("I128").Value = 170000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 510000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -505020
$ws.Range("H129").Value = 4167295.5
$ws.Range("I129").Value = 296.5
$ws.Range("J129").Value = 8334294.5
$ws.Range("K129").Value = 889.5
$ws.Range("L129").Value = 25002883.5
$ws.Range("M129").Value = 4110.5
$ws.Range("N129").Value = -25012883.5
$ws.Range("H130").Value = 2499.8572
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 2699.8
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 8099.400000000001
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -18139.4
$ws.Range("H131").Value = 20342.365
$ws.Range("I131").Value = 1203.2142
$ws.Range("J131").Value = 27393.63
$ws.Range("K131").Value = 3609.6426
$ws.Range("L131").Value = 82180.89
$ws.Range("M131").Value = 1430.3574
$ws.Range("N131").Value = -92260.89
$ws.Range("H132").Value = 1288.625
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1384.8334
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 12463.5006
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -17523.5006
$ws.Range("H133").Value = 4107.3477
$ws.Range("I133").Value = 2061
$ws.Range("J133").Value = 5681.4614
$ws.Range("K133").Value = 6183
$ws.Range("L133").Value = 17044.3842
$ws.Range("M133").Value = -1123
$ws.Range("N133").Value = -27164.3842
$ws.Range("H134").Value = 3197.4285
$ws.Range("I134").Value = 1757.1111
$ws.Range("J134").Value = 5790
$ws.Range("K134").Value = 5271.3333
$ws.Range("L134").Value = 17370
$ws.Range("M134").Value = -201.3333000000002
$ws.Range("N134").Value = -27510
$ws.Range("H136").Value = 3263.625
$ws.Range("I136").Value = 1165
$ws.Range("J136").Value = 3963.1667
$ws.Range("K136").Value = 3495
$ws.Range("L136").Value = 11889.5001
$ws.Range("M136").Value = 1605
$ws.Range("N136").Value = -22089.5001
$ws.Range("H137").Value = 55560304
$ws.Range("I137").Value = 4430
$ws.Range("J137").Value = 66671480
$ws.Range("K137").Value = 13290
$ws.Range("L137").Value = 200014440
$ws.Range("M137").Value = -8190
$ws.Range("N137").Value = -200024640
$ws.Range("H138").Value = 2918.2856
$ws.Range("I138").Value = 986
$ws.Range("J138").Value = 3991.7778
$ws.Range("K138").Value = 2958
$ws.Range("L138").Value = 11975.3334
$ws.Range("M138").Value = 2182
$ws.Range("N138").Value = -22255.3334
$ws.Range("H139").Value = 2256.8333
$ws.Range("I139").Value = 2057.7896
$ws.Range("J139").Value = 3013.2
$ws.Range("K139").Value = 6173.3688
$ws.Range("L139").Value = 9039.599999999999
$ws.Range("M139").Value = -1033.3688
$ws.Range("N139").Value = -19319.6
$ws.Range("H140").Value = 1473.9231
$ws.Range("I140").Value = 718.1667
$ws.Range("J140").Value = 3174.375
$ws.Range("K140").Value = 2154.5001
$ws.Range("L140").Value = 9523.125
$ws.Range("M140").Value = 3025.4999
$ws.Range("N140").Value = -19883.125
$ws.Range("H141").Value = 30021
$ws.Range("I141").Value = 20015
$ws.Range("J141").Value = 50033
$ws.Range("K141").Value = 60045
$ws.Range("L141").Value = 150099
$ws.Range("M141").Value = -54865
$ws.Range("N141").Value = -160459

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 18836
$ws.Range("I31").Value = 18836
$ws.Range("K31").Value = 18836
$ws.Range("M31").Value = -18544
$ws.Range("H37").Value = 18836
$ws.Range("I37").Value = 18836
$ws.Range("K37").Value = 18836
$ws.Range("M37").Value = -18559
$ws.Range("H134").Value = 92872.25
$ws.Range("J134").Value = 92872.25
$ws.Range("L134").Value = 278616.75
$ws.Range("N134").Value = -283686.75
$ws.Range("H136").Value = 27063
$ws.Range("J136").Value = 27063
$ws.Range("L136").Value = 81189
$ws.Range("N136").Value = -86289

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 731
$ws.Range("I35").Value = 731
$ws.Range("K35").Value = 731
$ws.Range("M35").Value = -395
$ws.Range("H132").Value = 4036.6487
$ws.Range("I132").Value = 3910.5293
$ws.Range("J132").Value = 5466
$ws.Range("K132").Value = 11731.5879
$ws.Range("L132").Value = 16398
$ws.Range("M132").Value = -9201.5879
$ws.Range("N132").Value = -21458
$ws.Range("H136").Value = 2969
$ws.Range("I136").Value = 2489.3333
$ws.Range("J136").Value = 3585.7144
$ws.Range("K136").Value = 7467.999899999999
$ws.Range("L136").Value = 10757.1432
$ws.Range("M136").Value = -4917.999899999999
$ws.Range("N136").Value = -15857.1432

Write-Host "Applied 248 cell updates across 7 worksheets."
